$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 21 for the "CRP_nr" summary line.
# This shifts the former rows 21 (crop_nr), 22 (forest_nr) and 23 (urban_nr)
# down to rows 22, 23 and 24 respectively, preserving their formatting.
$ws.Rows.Item(21).Insert()

# --- Updated summary statistics for the existing land-cover rows (2-20) ---
$ws.Range("B2").Value = 32.311275482177734
$ws.Range("C2").Value = 31.278425216674805
$ws.Range("D2").Value = 29.534662246704102
$ws.Range("E2").Value = 29.064609527587891
$ws.Range("F2").Value = 28.556692123413086
$ws.Range("G2").Value = 27.937778472900391
$ws.Range("H2").Value = 28.15052604675293

$ws.Range("C3").Value = 0.85814207792282104
$ws.Range("D3").Value = 2.1809043884277344
$ws.Range("E3").Value = 2.0747904777526855
$ws.Range("F3").Value = 1.8791482448577881
$ws.Range("G3").Value = 1.9326847791671753
$ws.Range("H3").Value = 1.4011971950531006

$ws.Range("B4").Value = 34.885387420654297
$ws.Range("C4").Value = 35.07537841796875
$ws.Range("D4").Value = 35.097221374511719
$ws.Range("E4").Value = 35.111167907714844
$ws.Range("F4").Value = 35.214580535888672
$ws.Range("G4").Value = 35.234233856201172
$ws.Range("H4").Value = 35.320003509521484

$ws.Range("B5").Value = 11.892159461975098
$ws.Range("C5").Value = 11.488377571105957
$ws.Range("D5").Value = 11.323563575744629
$ws.Range("E5").Value = 10.844473838806152
$ws.Range("F5").Value = 10.637063026428223
$ws.Range("G5").Value = 10.674258232116699
$ws.Range("H5").Value = 10.687778472900391

$ws.Range("B6").Value = 15.533368110656738
$ws.Range("C6").Value = 15.299148559570313
$ws.Range("D6").Value = 15.149826049804688
$ws.Range("E6").Value = 15.119749069213867
$ws.Range("F6").Value = 15.110164642333984
$ws.Range("G6").Value = 15.08366870880127
$ws.Range("H6").Value = 15.032883644104004

$ws.Range("B7").Value = 5.3127474784851074
$ws.Range("C7").Value = 5.9354662895202637
$ws.Range("D7").Value = 6.6487612724304199
$ws.Range("E7").Value = 7.7201461791992188
$ws.Range("F7").Value = 8.5372905731201172
$ws.Range("G7").Value = 9.0723123550415039
$ws.Range("H7").Value = 9.3425521850585938

$ws.Range("B8").Value = 17.206083297729492
$ws.Range("C8").Value = 17.7904052734375
$ws.Range("D8").Value = 18.486202239990234
$ws.Range("E8").Value = 19.427421569824219
$ws.Range("F8").Value = 20.200677871704102
$ws.Range("G8").Value = 20.673070907592773
$ws.Range("H8").Value = 20.925247192382813

$ws.Range("B9").Value = 2.1980595588684082
$ws.Range("C9").Value = 2.1871206760406494
$ws.Range("D9").Value = 2.1753609180450439
$ws.Range("E9").Value = 2.1508204936981201
$ws.Range("F9").Value = 2.1286287307739258
$ws.Range("G9").Value = 2.1121630668640137
$ws.Range("H9").Value = 2.1029939651489258

$ws.Range("B10").Value = 21.829973220825195
$ws.Range("C10").Value = 21.649463653564453
$ws.Range("D10").Value = 21.448871612548828
$ws.Range("E10").Value = 21.175046920776367
$ws.Range("F10").Value = 20.949491500854492
$ws.Range("G10").Value = 20.799098968505859
$ws.Range("H10").Value = 20.726394653320313

$ws.Range("B11").Value = 18.786640167236328
$ws.Range("C11").Value = 18.659990310668945
$ws.Range("D11").Value = 18.494180679321289
$ws.Range("E11").Value = 18.266914367675781
$ws.Range("F11").Value = 18.096242904663086
$ws.Range("G11").Value = 17.997442245483398
$ws.Range("H11").Value = 17.940835952758789

$ws.Range("B12").Value = 11.753081321716309
$ws.Range("C12").Value = 11.663716316223145
$ws.Range("D12").Value = 11.548404693603516
$ws.Range("E12").Value = 11.393607139587402
$ws.Range("F12").Value = 11.266033172607422
$ws.Range("G12").Value = 11.186065673828125
$ws.Range("H12").Value = 11.145540237426758

$ws.Range("B13").Value = 2.373673677444458
$ws.Range("C13").Value = 2.3585808277130127
$ws.Range("D13").Value = 2.3385131359100342
$ws.Range("E13").Value = 2.3175194263458252
$ws.Range("F13").Value = 2.3024177551269531
$ws.Range("G13").Value = 2.2910037040710449
$ws.Range("H13").Value = 2.2844910621643066

$ws.Range("B14").Value = 12.506346702575684
$ws.Range("C14").Value = 12.433952331542969
$ws.Range("D14").Value = 12.348419189453125
$ws.Range("E14").Value = 12.227241516113281
$ws.Range("F14").Value = 12.130745887756348
$ws.Range("G14").Value = 12.07868480682373
$ws.Range("H14").Value = 12.046497344970703

$ws.Range("B15").Value = 12.40440845489502
$ws.Range("C15").Value = 12.326560974121094
$ws.Range("D15").Value = 12.241860389709473
$ws.Range("E15").Value = 12.138484001159668
$ws.Range("F15").Value = 12.028786659240723
$ws.Range("G15").Value = 11.966737747192383
$ws.Range("H15").Value = 11.929265022277832

$ws.Range("B16").Value = 0.94173288345336914
$ws.Range("C16").Value = 0.9302094578742981
$ws.Range("D16").Value = 0.91818696260452271
$ws.Range("E16").Value = 0.90294367074966431
$ws.Range("F16").Value = 0.89697593450546265
$ws.Range("G16").Value = 0.89573168754577637
$ws.Range("H16").Value = 0.89873361587524414

$ws.Range("B17").Value = 24.028032302856445
$ws.Range("C17").Value = 23.836584091186523
$ws.Range("D17").Value = 23.624233245849609
$ws.Range("E17").Value = 23.32586669921875
$ws.Range("F17").Value = 23.078119277954102
$ws.Range("G17").Value = 22.911262512207031
$ws.Range("H17").Value = 22.829387664794922

$ws.Range("B18").Value = 30.539722442626953
$ws.Range("C18").Value = 30.323707580566406
$ws.Range("D18").Value = 30.042585372924805
$ws.Range("E18").Value = 29.660520553588867
$ws.Range("F18").Value = 29.362276077270508
$ws.Range("G18").Value = 29.183507919311523
$ws.Range("H18").Value = 29.086376190185547

$ws.Range("B19").Value = 14.880021095275879
$ws.Range("C19").Value = 14.792532920837402
$ws.Range("D19").Value = 14.686932563781738
$ws.Range("E19").Value = 14.544760704040527
$ws.Range("F19").Value = 14.433163642883301
$ws.Range("G19").Value = 14.369688034057617
$ws.Range("H19").Value = 14.330988883972168

$ws.Range("B20").Value = 13.346141815185547
$ws.Range("C20").Value = 13.256770133972168
$ws.Range("D20").Value = 13.16004753112793
$ws.Range("E20").Value = 13.041427612304688
$ws.Range("F20").Value = 12.925762176513672
$ws.Range("G20").Value = 12.862469673156738
$ws.Range("H20").Value = 12.827999114990234

# --- New row 21: CRP_nr (B21 intentionally left blank) ---
$ws.Range("A21").Value = "CRP_nr"
$ws.Range("C21").Value = 50.865180969238281
$ws.Range("D21").Value = 52.805267333984375
$ws.Range("E21").Value = 53.657325744628906
$ws.Range("F21").Value = 55.877994537353516
$ws.Range("G21").Value = 59.601699829101563
$ws.Range("H21").Value = 68.188255310058594

# --- Row 23 (was row 22, "forest_nr"): refreshed F/G/H figures ---
$ws.Range("F23").Value = 20.555309295654297
$ws.Range("G23").Value = 19.22590446472168
$ws.Range("H23").Value = 16.138092041015625

# --- Row 24 (was row 23, "urban_nr"): refreshed E/F/G/H figures ---
$ws.Range("E24").Value = 22334.931640625
$ws.Range("F24").Value = 28291.578125
$ws.Range("G24").Value = 41637.09375
$ws.Range("H24").Value = 43642.6015625
